$d = $word.ActiveDocument

# --- Hunk 1: parcels ownershipType placeholder (drop ".label") ---
$d.Content.Find.Execute(
    "{d.parcels[i].ownershipType.label:ifEM():show(.noData)}", $true, $false, $false, $false, $false,
    $true, 1, $false, "{d.parcels[i].ownershipType:ifEM():show(.noData)}", 2
) | Out-Null

# --- Hunk 3: otherParcels ownershipType placeholder (drop ".label") ---
$d.Content.Find.Execute(
    "{d.otherParcels[i].ownershipType.label:ifEM():show(.noData)}", $true, $false, $false, $false, $false,
    $true, 1, $false, "{d.otherParcels[i].ownershipType:ifEM():show(.noData)}", 2
) | Out-Null

# --- Hunk 2 & 4: the empty "Owner information" spacer paragraphs (Heading 3,
# no text) pick up the same run formatting as the paragraph mark right next
# to them (gray 565656, 13pt, not bold) instead of just the BC Sans font
# override. Re-apply that formatting to the lone empty run in each paragraph.
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text
    $styleName = $p.Style.NameLocal
    if ($styleName -eq "Heading 3" -and $txt.Length -le 1) {
        $r = $p.Range
        $r.Font.Bold = $false
        $r.Font.BoldBi = $false
        $r.Font.Color = 5658198
        $r.Font.Size = 13
        $r.Font.SizeBi = 13
    }
}
